# Apply the cell-level changes described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set new/changed cell values (order chosen so we never overwrite a value
# before it has been read/moved elsewhere).
$ws.Range("B5").Value = "ActivatedItem"
$ws.Range("B6").Value = "ActivatedItemMGr"

$ws.Range("A8").Value = "AnimatSim::Behavior"
$ws.Range("B9").Value = "NeuralModule"
$ws.Range("B10").Value = "NervousSystem"

$ws.Range("A12").Value = "AnimatSim::Environment"
$ws.Range("B13").Value = "Structure"
$ws.Range("B14").Value = "organism"

# Clear the now-vacated cells from the old layout.
$ws.Range("A6").Value = $null
$ws.Range("B7").Value = $null
$ws.Range("B8").Value = $null
$ws.Range("A10").Value = $null
$ws.Range("B11").Value = $null
$ws.Range("B12").Value = $null

# Update the selection to match the saved view state.
$ws.Range("A12").Select()
